# Build site at 2022-09-26 16:07:08 UTC
#
# The source "Docentes responsáveis:" value ("11079086 - Herlandí de Souza
# Andrade") used to live in its own row (13) underneath the "Docentes
# responsáveis:" label row (12). The edit removes that standalone row and
# re-purposes the text elsewhere in the sheet (it now shows up as the
# "Objetivos:" value and, further down, as the "Método:" value), while the
# remaining rows shift up by one and several cell contents below get
# shuffled around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old standalone row that only held the docente name
#    (row 13, with empty column A) -- everything below shifts up by one.
$ws.Rows.Item(13).Delete()

# 2) Now patch the handful of cells whose text content changed as part of
#    the edit (row numbers below are POST-shift, i.e. matching the final
#    layout).

# "Objetivos:" value becomes the docente name text.
$ws.Cells.Item(10, 2).Value2 = "11079086 - Herlandí de Souza Andrade"
$ws.Cells.Item(10, 3).Value2 = "11079086 - Herlandí de Souza Andrade"

# "Programa resumido:" value becomes "Semestral".
$ws.Cells.Item(13, 2).Value2 = "Semestral"
$ws.Cells.Item(13, 3).Value2 = "Semestral"

# "Programa:" value becomes "01/01/2021". Copy/paste-special (values only)
# from the existing "Ativação:" value cells so the text lands as a plain
# shared-string cell instead of Excel auto-parsing the dd/mm/yyyy-looking
# text into a date serial number (which would also mint a new style).
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# "Método:" value becomes the docente name text.
$ws.Cells.Item(18, 2).Value2 = "11079086 - Herlandí de Souza Andrade"
$ws.Cells.Item(18, 3).Value2 = "11079086 - Herlandí de Souza Andrade"

# "Critério:" value becomes the A1/A2 evaluation paragraph.
$ws.Cells.Item(19, 2).Value2 = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Cells.Item(19, 3).Value2 = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."

# "Norma de recuperação:" value becomes the weighted-average sentence.
$ws.Cells.Item(20, 2).Value2 = "Média ponderada das avaliações (M)."
$ws.Cells.Item(20, 3).Value2 = "Média ponderada das avaliações (M)."

# "Bibliografia:" value becomes the recovery-exam paragraph.
$ws.Cells.Item(21, 2).Value2 = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Cells.Item(21, 3).Value2 = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
